# Change `NilaiPeriodik` from @Embeddable to entity:
#  - remove the old "periodekas_listtransaksikas" embedded-collection sheet
#  - add "periodekas_transaksikas" (lean join table: PeriodeKas_id / listTransaksiKas_id / listTransaksiKas_ORDER)
#  - add "transaksikas" (the promoted entity, with full entity columns)
#
# New sheet content is written FIRST (while the old sheet is still present)
# so every shared string the old sheet used stays "in use" and keeps its
# original index; only the brand-new string ("listTransaksiKas_id") gets
# appended at the end. The old sheet is deleted last.

$wb = $excel.ActiveWorkbook

$oldWs = $wb.Worksheets.Item("periodekas_listtransaksikas")
$lastWs = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---- new sheet: periodekas_transaksikas --------------------------------
$joinWs = $wb.Worksheets.Add($null, $lastWs)
$joinWs.Name = "periodekas_transaksikas"

$joinWs.Range("A1").Value = "PeriodeKas_id"
$joinWs.Range("B1").Value = "listTransaksiKas_id"
$joinWs.Range("C1").Value = "listTransaksiKas_ORDER"

$joinWs.Range("A2").Value = -1
$joinWs.Range("B2").Value = -1
$joinWs.Range("C2").Value = 0

$joinWs.Range("A3").Value = -1
$joinWs.Range("B3").Value = -2
$joinWs.Range("C3").Value = 1

# header row styling (bold / grey-filled header style already used elsewhere)
$wb.Worksheets.Item("kas").Range("A1").Copy() | Out-Null
$joinWs.Range("A1:C1").PasteSpecial(-4122) | Out-Null

$joinWs.Range("C2").Select() | Out-Null

# ---- new sheet: transaksikas (promoted entity) --------------------------
$txWs = $wb.Worksheets.Add($null, $joinWs)
$txWs.Name = "transaksikas"

$txWs.Range("A1").Value = "id"
$txWs.Range("B1").Value = "jumlah"
$txWs.Range("C1").Value = "keterangan"
$txWs.Range("D1").Value = "saldo"
$txWs.Range("E1").Value = "tanggal"
$txWs.Range("F1").Value = "createdBy"
$txWs.Range("G1").Value = "createdDate"
$txWs.Range("H1").Value = "deleted"
$txWs.Range("I1").Value = "modifiedBy"
$txWs.Range("J1").Value = "modifiedDate"
$txWs.Range("K1").Value = "pihakTerkait"
$txWs.Range("L1").Value = "jenis_id"
$txWs.Range("M1").Value = "kategoriKas_id"

$txWs.Range("A2").Value = -1
$txWs.Range("B2").Value = 10000
$txWs.Range("D2").Value = 10000
$txWs.Range("K2").Value = "unknown"
$txWs.Range("L2").Value = -1
$txWs.Range("M2").Value = -1

$txWs.Range("A3").Value = -2
$txWs.Range("B3").Value = 12000
$txWs.Range("D3").Value = 22000
$txWs.Range("K3").Value = "snake"
$txWs.Range("L3").Value = -2
$txWs.Range("M3").Value = -1

# header style
$wb.Worksheets.Item("kas").Range("A1").Copy() | Out-Null
$txWs.Range("A1:M1").PasteSpecial(-4122) | Out-Null

# date cells (E2/E3) — copy date style+value from the source sheet being retired
# (value first, then the format-only paste, so the paste doesn't get clobbered)
$txWs.Range("E2").Value = 41641
$txWs.Range("E3").Value = 41642
$oldWs.Range("H2").Copy() | Out-Null
$txWs.Range("E2").PasteSpecial(-4122) | Out-Null
$oldWs.Range("H3").Copy() | Out-Null
$txWs.Range("E3").PasteSpecial(-4122) | Out-Null

# pihakTerkait cells (K2/K3) — copy the "admin"-style text formatting
# (value first, then the format-only paste, so the paste doesn't get clobbered)
$txWs.Range("K2").Value = "unknown"
$txWs.Range("K3").Value = "snake"
$oldWs.Range("D2").Copy() | Out-Null
$txWs.Range("K2").PasteSpecial(-4122) | Out-Null
$oldWs.Range("D3").Copy() | Out-Null
$txWs.Range("K3").PasteSpecial(-4122) | Out-Null

$txWs.Range("M3").Select() | Out-Null

# ---- drop the old embedded-collection sheet ------------------------------
$oldWs.Delete()

# re-activate periodekas_transaksikas so it carries tabSelected="1"
# (re-fetch by name: the old $joinWs handle goes stale across the Delete())
$joinWsAfter = $wb.Worksheets.Item("periodekas_transaksikas")
$joinWsAfter.Activate()
$joinWsAfter.Range("C2").Select() | Out-Null
